# Weekly update: a new price record was reported for the week of 2021-09-29
# (serial date 44468). Insert it above the existing row 84 so the whole
# table shifts down by one row (old row 84 becomes row 85, ..., old row 116
# becomes row 117) and fill in the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 84:116 down to 85:117, carrying formatting with them.
$ws.Rows("84:84").Insert()

# Populate the newly-inserted row 84 with the new weekly record.
$ws.Range("A84").Value = 4
$ws.Range("B84").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C84").Value = "Los Lagos"
$ws.Range("D84").Value = 44468
$ws.Range("E84").Value = 10
$ws.Range("F84").Value = "Fruta"
$ws.Range("G84").Value = 100101
$ws.Range("H84").Value = "Berries"
$ws.Range("I84").Value = 100101007
$ws.Range("J84").Value = "Kiwi"
$ws.Range("K84").Value = "Hayward"
$ws.Range("L84").Value = "Primera"
$ws.Range("M84").Value = 150
$ws.Range("N84").Value = 15000
$ws.Range("O84").Value = 15000
$ws.Range("P84").Value = 15000
$ws.Range("Q84").Value = "$/caja 15 kilos"
$ws.Range("R84").Value = "Provincia de Curicó"
$ws.Range("S84").Value = 1000
$ws.Range("T84").Value = 15
